$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = 482
$ws.Range("C2:C$lastRow").Value = 45179
